$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 48) with the latest ranking snapshot.
$ws.Range("A48").Value = "2025/12/04 17:00"
$ws.Range("B48").Value = "-"
$ws.Range("C48").Value = "-"
$ws.Range("D48").Value = "-"
$ws.Range("E48").Value = "-"
$ws.Range("F48").Value = "-"
$ws.Range("G48").Value = "-"
